$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 735.3077
$ws.Range("I18").Value = 383.625
$ws.Range("K18").Value = 383.625
$ws.Range("M18").Value = -99.625

$ws.Range("H19").Value = 546.2
$ws.Range("I19").Value = 642.8182
$ws.Range("K19").Value = 642.8182
$ws.Range("M19").Value = -467.8182

$ws.Range("H32").Value = 400
$ws.Range("J32").Value = 400
$ws.Range("L32").Value = 400
$ws.Range("N32").Value = -1052

$ws.Range("H33").Value = 601.15
$ws.Range("I33").Value = 643.94446
$ws.Range("K33").Value = 643.94446
$ws.Range("M33").Value = -414.94446

$ws.Range("H40").Value = 2003.7273
$ws.Range("I40").Value = 1793.6364
$ws.Range("J40").Value = 2213.818
$ws.Range("K40").Value = 1793.6364
$ws.Range("L40").Value = 2213.818
$ws.Range("M40").Value = -1618.6364
$ws.Range("N40").Value = -2563.818

$ws.Range("H76").Value = 151254
$ws.Range("I76").Value = 372870.34
$ws.Range("J76").Value = 3509.7778
$ws.Range("K76").Value = 372870.34
$ws.Range("L76").Value = 3509.7778
$ws.Range("M76").Value = -372555.34
$ws.Range("N76").Value = -4139.7778

$ws.Range("H79").Value = 151254
$ws.Range("I79").Value = 372870.34
$ws.Range("J79").Value = 3509.7778
$ws.Range("K79").Value = 372870.34
$ws.Range("L79").Value = 3509.7778
$ws.Range("M79").Value = -371778.34
$ws.Range("N79").Value = -5693.7778

$ws.Range("H112").Value = 1417.3334
$ws.Range("I112").Value = 800
$ws.Range("J112").Value = 1494.5
$ws.Range("K112").Value = 2400
$ws.Range("L112").Value = 4483.5
$ws.Range("M112").Value = -1292
$ws.Range("N112").Value = -6699.5

$ws.Range("H113").Value = 3268.75
$ws.Range("I113").Value = 2720
$ws.Range("J113").Value = 3518.182
$ws.Range("K113").Value = 2720
$ws.Range("L113").Value = 3518.182
$ws.Range("M113").Value = 534
$ws.Range("N113").Value = -10026.182

$ws.Range("H127").Value = 872.5
$ws.Range("I127").Value = 619.3333
$ws.Range("J127").Value = 1198
$ws.Range("K127").Value = 1857.9999
$ws.Range("L127").Value = 3594
$ws.Range("M127").Value = 3102.0001
$ws.Range("N127").Value = -13514

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 784.3077
$ws.Range("I2").Value = 770.5454999999999
$ws.Range("K2").Value = 770.5454999999999
$ws.Range("M2").Value = -657.5454999999999

$ws.Range("H50").Value = 733.3333
$ws.Range("I50").Value = 360
$ws.Range("J50").Value = 1032
$ws.Range("K50").Value = 360
$ws.Range("L50").Value = 1032
$ws.Range("M50").Value = 354
$ws.Range("N50").Value = -2460

$ws.Range("H63").Value = 4357.857
$ws.Range("I63").Value = 2401.25
$ws.Range("J63").Value = 6966.6665
$ws.Range("K63").Value = 2401.25
$ws.Range("L63").Value = 6966.6665
$ws.Range("M63").Value = -1715.25
$ws.Range("N63").Value = -8338.666499999999

$ws.Range("H66").Value = 4357.857
$ws.Range("I66").Value = 2401.25
$ws.Range("J66").Value = 6966.6665
$ws.Range("K66").Value = 12006.25
$ws.Range("L66").Value = 34833.3325
$ws.Range("M66").Value = -8574.25
$ws.Range("N66").Value = -41697.3325

$ws.Range("H116").Value = 784.3077
$ws.Range("I116").Value = 770.5454999999999
$ws.Range("K116").Value = 770.5454999999999
$ws.Range("M116").Value = 1523.4545

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 784.3077
$ws.Range("I3").Value = 770.5454999999999
$ws.Range("K3").Value = 770.5454999999999
$ws.Range("M3").Value = -656.5454999999999

$ws.Range("H15").Value = 3336
$ws.Range("J15").Value = 3336
$ws.Range("L15").Value = 3336
$ws.Range("N15").Value = -3790

$ws.Range("H20").Value = 2162.3809
$ws.Range("I20").Value = 2163
$ws.Range("J20").Value = 2162
$ws.Range("K20").Value = 2163
$ws.Range("L20").Value = 2162
$ws.Range("M20").Value = -1916
$ws.Range("N20").Value = -2656

$ws.Range("H134").Value = 2772.5
$ws.Range("I134").Value = 2296.923
$ws.Range("K134").Value = 6890.768999999999
$ws.Range("M134").Value = -4355.768999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 716.6667
$ws.Range("I22").Value = 625
$ws.Range("J22").Value = 900
$ws.Range("K22").Value = 625
$ws.Range("L22").Value = 900
$ws.Range("M22").Value = -275
$ws.Range("N22").Value = -1600

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1667659.2
$ws.Range("I131").Value = 11111534
$ws.Range("J131").Value = 1093.0883
$ws.Range("K131").Value = 33334602
$ws.Range("L131").Value = 3279.2649
$ws.Range("M131").Value = -33329562
$ws.Range("N131").Value = -13359.2649

$ws.Range("H140").Value = 1466.6333
$ws.Range("I140").Value = 994.3889
$ws.Range("K140").Value = 2983.1667
$ws.Range("M140").Value = 2196.8333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 972.381
$ws.Range("I16").Value = 648
$ws.Range("J16").Value = 1783.3334
$ws.Range("K16").Value = 648
$ws.Range("L16").Value = 1783.3334
$ws.Range("M16").Value = -478
$ws.Range("N16").Value = -2123.3334

$ws.Range("H22").Value = 832.25
$ws.Range("I22").Value = 427
$ws.Range("J22").Value = 1237.5
$ws.Range("K22").Value = 427
$ws.Range("L22").Value = 1237.5
$ws.Range("M22").Value = -132
$ws.Range("N22").Value = -1827.5

$ws.Range("H27").Value = 832.25
$ws.Range("I27").Value = 427
$ws.Range("J27").Value = 1237.5
$ws.Range("K27").Value = 427
$ws.Range("L27").Value = 1237.5
$ws.Range("M27").Value = -320
$ws.Range("N27").Value = -1451.5

$ws.Range("H46").Value = 1023.1539
$ws.Range("I46").Value = 400.5
$ws.Range("J46").Value = 1136.3636
$ws.Range("K46").Value = 400.5
$ws.Range("L46").Value = 1136.3636
$ws.Range("M46").Value = -212.5
$ws.Range("N46").Value = -1512.3636

$ws.Range("H82").Value = 1088.04
$ws.Range("I82").Value = 875.9167
$ws.Range("J82").Value = 1283.8462
$ws.Range("K82").Value = 875.9167
$ws.Range("L82").Value = 1283.8462
$ws.Range("M82").Value = -514.9167
$ws.Range("N82").Value = -2005.8462

$ws.Range("H85").Value = 1088.04
$ws.Range("I85").Value = 875.9167
$ws.Range("J85").Value = 1283.8462
$ws.Range("K85").Value = 875.9167
$ws.Range("L85").Value = 1283.8462
$ws.Range("M85").Value = 372.0833
$ws.Range("N85").Value = -3779.8462

$ws.Range("H104").Value = 8274
$ws.Range("J104").Value = 8274
$ws.Range("L104").Value = 8274
$ws.Range("N104").Value = -15262

$ws.Range("H132").Value = 8626471
$ws.Range("I132").Value = 19240828
$ws.Range("J132").Value = 2306.4688
$ws.Range("K132").Value = 57722484
$ws.Range("L132").Value = 6919.4064
$ws.Range("M132").Value = -57719954
$ws.Range("N132").Value = -11979.4064

$ws.Range("H140").Value = 30000
$ws.Range("J140").Value = 30000
$ws.Range("L140").Value = 30000
$ws.Range("N140").Value = -40360

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H138").Value = 80878.625
$ws.Range("J138").Value = 80878.625
$ws.Range("L138").Value = 80878.625
$ws.Range("N138").Value = -91158.625
